# Apply the Natmi (Il16-Grin2d) data refresh described in the commit
# "Natmi following Dr Hou advice": recompute the interaction table with
# updated cell counts / expression statistics, and extend the table from
# 12 data rows (4 sending clusters x 3 target clusters) to 16 data rows
# (4 sending clusters x 4 target clusters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("ECs","Il16","Grin2d","ECs",3,1,4.243623333333334,12.73087,0.2469246453968972,0.2469246453968973,3,1,0.1972836666666667,0.591851,0.08392655655245007,0.08392655655245007,0.8371975711522223,7.534778140370001,0.02072353521609637,0.02072353521609638)
  ,@("ECs","Il16","Grin2d","FAPs",3,1,4.243623333333334,12.73087,0.2469246453968972,0.2469246453968973,3,1,1.922001,5.766003,0.8176395357295954,0.8176395357295954,8.15624829029,73.40623461261,0.201895352422514,0.2018953524225141)
  ,@("ECs","Il16","Grin2d","M2",3,1,4.243623333333334,12.73087,0.2469246453968972,0.2469246453968973,2,0.6666666666666666,0.03285833333333334,0.09857500000000001,0.01397828222332608,0.01397828222332608,0.1394383900277778,1.25494551025,0.003451582381252546,0.003451582381252546)
  ,@("FAPs","Il16","Grin2d","sCs",3,1,4.243623333333334,12.73087,0.2469246453968972,0.2469246453968973,3,1,0.1985273333333334,0.5955820000000001,0.08445562549462841,0.08445562549462841,0.842475224037778,7.582277016340002,0.02085417537703427,0.02085417537703428)
  ,@("FAPs","Il16","Grin2d","ECs",3,1,5.865491666666667,17.596475,0.3412966552647515,0.3412966552647516,3,1,0.1972836666666667,0.591851,0.08392655655245007,0.08392655655245007,1.157165702802778,10.414491325225,0.02864385303923922,0.02864385303923923)
  ,@("FAPs","Il16","Grin2d","FAPs",3,1,5.865491666666667,17.596475,0.3412966552647515,0.3412966552647516,3,1,1.922001,5.766003,0.8176395357295954,0.8176395357295954,11.273480848825,101.461327639425,0.2790576387567352,0.2790576387567353)
  ,@("FAPs","Il16","Grin2d","M2",3,1,5.865491666666667,17.596475,0.3412966552647515,0.3412966552647516,2,0.6666666666666666,0.03285833333333334,0.09857500000000001,0.01397828222332608,0.01397828222332608,0.1927302803472223,1.734572523125,0.004770740969167927,0.004770740969167929)
  ,@("FAPs","Il16","Grin2d","sCs",3,1,5.865491666666667,17.596475,0.3412966552647515,0.3412966552647516,3,1,0.1985273333333334,0.5955820000000001,0.08445562549462841,0.08445562549462841,1.164460419272223,10.48014377345,0.02882442249960915,0.02882442249960916)
  ,@("M2","Il16","Grin2d","ECs",3,1,6.123111999999999,18.369336,0.3562868663317164,0.3562868663317164,3,1,0.1972836666666667,0.591851,0.08392655655245007,0.08392655655245007,1.207989986770666,10.871909880936,0.02990192983608401,0.02990192983608402)
  ,@("M2","Il16","Grin2d","FAPs",3,1,6.123111999999999,18.369336,0.3562868663317164,0.3562868663317164,3,1,1.922001,5.766003,0.8176395357295954,0.8176395357295954,11.768627387112,105.917646484008,0.291314227974017,0.291314227974017)
  ,@("M2","Il16","Grin2d","M2",3,1,6.123111999999999,18.369336,0.3562868663317164,0.3562868663317164,2,0.6666666666666666,0.03285833333333334,0.09857500000000001,0.01397828222332608,0.01397828222332608,0.2011952551333333,1.8107572962,0.004980278370049188,0.004980278370049189)
  ,@("M2","Il16","Grin2d","sCs",3,1,6.123111999999999,18.369336,0.3562868663317164,0.3562868663317164,3,1,0.1985273333333334,0.5955820000000001,0.08445562549462841,0.08445562549462841,1.215605097061333,10.940445873552,0.03009043015156617,0.03009043015156617)
  ,@("sCs","Il16","Grin2d","ECs",3,1,0.9536773333333334,2.861032,0.05549183300663471,0.05549183300663472,3,1,0.1972836666666667,0.591851,0.08392655655245007,0.08392655655245007,0.1881449611368889,1.693304650232,0.004657238461030443,0.004657238461030445)
  ,@("sCs","Il16","Grin2d","FAPs",3,1,0.9536773333333334,2.861032,0.05549183300663471,0.05549183300663472,3,1,1.922001,5.766003,0.8176395357295954,0.8176395357295954,1.832968788344,16.496719095096,0.04537231657632904,0.04537231657632905)
  ,@("sCs","Il16","Grin2d","M2",3,1,0.9536773333333334,2.861032,0.05549183300663471,0.05549183300663472,2,0.6666666666666666,0.03285833333333334,0.09857500000000001,0.01397828222332608,0.01397828222332608,0.03133624771111112,0.2820262294000001,0.0007756805028564217,0.0007756805028564218)
  ,@("sCs","Il16","Grin2d","sCs",3,1,0.9536773333333334,2.861032,0.05549183300663471,0.05549183300663472,3,1,0.1985273333333334,0.5955820000000001,0.08445562549462841,0.08445562549462841,0.1893310178471111,1.703979160624,0.004686597466418801,0.004686597466418802)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $rowNum = $i + 2
  $rowData = $data[$i]
  for ($j = 0; $j -lt $rowData.Count; $j++) {
    $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
  }
}
